$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cols = $ws.Range($ws.Columns.Item(36), $ws.Columns.Item(39))
$cols.EntireColumn.Delete(-4159)
